$wb = $excel.ActiveWorkbook

$wsGTGT = $wb.Worksheets.Item("GTGT")
$wsTNDN = $wb.Worksheets.Item("TNDN")
$wsMonBai = $wb.Worksheets.Item("MÔN BÀI")

# ---------------------------------------------------------------------------
# Sheet "GTGT": fix C2 number format (was missing the #,##0 style) and fill in
# the 2020 quarterly GTGT tax figures (rows 19-22 of the 3rd mini-table).
# ---------------------------------------------------------------------------
$wsGTGT.Range("C2").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'

$wsGTGT.Range("B19").Value = 4319651
$wsGTGT.Range("C19").Value = 4319651
$wsGTGT.Range("E19").ClearContents()

$wsGTGT.Range("B20").Value = 3722955
$wsGTGT.Range("C20").Value = 3722955
$wsGTGT.Range("E20").ClearContents()

$wsGTGT.Range("B21").Value = 8075555
$wsGTGT.Range("C21").Value = 8075555
$wsGTGT.Range("E21").ClearContents()

$wsGTGT.Range("B22").Value = 3989945

# ---------------------------------------------------------------------------
# Sheet "TNDN": fill in the 2020 quarterly TNDN tax figures (rows 19-22 of
# the 3rd mini-table).
# ---------------------------------------------------------------------------
$wsTNDN.Range("B19").Value = 7354211.6000000006
$wsTNDN.Range("C19").Value = 1470842

$wsTNDN.Range("B20").Value = 6684178.2000000002
$wsTNDN.Range("C20").Value = 5000000

$wsTNDN.Range("B21").Value = 14518739.085601091
$wsTNDN.Range("C21").Value = 12000000

$wsTNDN.Range("B22").Value = 7164720.5594954491

# ---------------------------------------------------------------------------
# Sheet "MÔN BÀI": add the 2021 license-tax mini-table (header + data row +
# total row), mirroring the layout already used for 2018/2019/2020.
# ---------------------------------------------------------------------------
$wsMonBai.Range("A2:E2").Copy()
$wsMonBai.Range("A14:E14").PasteSpecial(-4122)
$wsMonBai.Range("A14").Value = "KỲ THUẾ"
$wsMonBai.Range("B14").Value = "KHAI THUẾ"
$wsMonBai.Range("C14").Value = "ĐÓNG THUẾ"
$wsMonBai.Range("D14").Value = "SỐ THUẾ CÒN PHẢI ĐÓNG "
$wsMonBai.Range("E14").Value = "NGÀY ĐÓNG THUẾ"

$wsMonBai.Range("A3:E3").Copy()
$wsMonBai.Range("A15:E15").PasteSpecial(-4122)
$wsMonBai.Range("A15").Value = "THUẾ MÔN BÀI NĂM 2021"
$wsMonBai.Range("B15").Value = 2000000
$wsMonBai.Range("C15").Value = 2000000
$wsMonBai.Range("D15").Formula = "=B15-C15"

$wsMonBai.Range("A4:E4").Copy()
$wsMonBai.Range("A16:E16").PasteSpecial(-4122)
$wsMonBai.Range("A16").Value = "TỔNG"
$wsMonBai.Range("B16").Formula = "=SUM(B15:B15)"
$wsMonBai.Range("C16").Formula = "=SUM(C15:C15)"
$wsMonBai.Range("D16").Formula = "=SUM(D15:D15)"

# ---------------------------------------------------------------------------
# Selection / active-cell bookkeeping to match the saved UI state. "GTGT"
# stays the active tab, so it must be the last sheet whose range is selected
# without an explicit Activate() on the other sheets.
# ---------------------------------------------------------------------------
$wsTNDN.Range("C22").Select()
$wsMonBai.Range("A19").Select()
$wsGTGT.Activate()
$wsGTGT.Range("C22").Select()
